$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Correct the slightly-off timestamp on row 15 (A15)
$ws.Cells.Item(15, 1).Value = 45868.66690025463

# Append the new automated-task row (row 16)
$ws.Cells.Item(16, 1).Value = 45868.70855332787
$ws.Cells.Item(16, 2).Value = 2025
$ws.Cells.Item(16, 3).Value = 31
$ws.Cells.Item(16, 4).Value = 18.72
$ws.Cells.Item(16, 5).Value = 78.18000000000001
$ws.Cells.Item(16, 6).Value = 114.26
$ws.Cells.Item(16, 7).Value = 15.74
$ws.Cells.Item(16, 8).Value = "ESE"
$ws.Cells.Item(16, 9).Value = 0
$ws.Cells.Item(16, 10).Value = "17:00:19"

# Match the existing date/time number-format style used by column A
$ws.Cells.Item(16, 1).NumberFormat = $ws.Cells.Item(15, 1).NumberFormat
